$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.551.79"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +5.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.508.98"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +8.39%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +10.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +10.92%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.501.81"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +8.38%  "

# Row 9
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.643"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +9.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +19.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.64"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.86%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000278"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +10.98%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.43"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.070.25"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +8.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.508.27"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.65"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +9.86%  "

# Row 18
$ws.Range("E18").Value = "  +6.53%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.587.42"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.79%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +11.19%  "

# Row 21
$ws.Range("E21").Value = "  +7.43%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +13.85%  "

# Row 23
$ws.Range("E23").Value = "  +13.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.99"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.66%  "

# Row 27
$ws.Range("E27").Value = "  +11.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.39"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +12.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.23"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +16.81%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.09%  "

# Row 32
$ws.Range("E32").Value = "  +3.89%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "619.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.93"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +8.93%  "

# Row 35
$ws.Range("E35").Value = "  +9.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.02"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +24.02%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0817"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +16.11%  "

# Row 39
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.32"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +9.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.09%  "

# Row 41
$ws.Range("E41").Value = "  +5.18%  "

# Row 42
$ws.Range("E42").Value = "  +10.83%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.127.78"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +11.33%  "

# Row 44
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.63"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.88%  "

# Row 46
$ws.Range("E46").Value = "  +14.05%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +12.29%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0418"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.06%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.74"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.33%  "

# Row 50
$ws.Range("E50").Value = "  +9.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.69"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.57%  "
